$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: Date, Model, Description, Score
# Copy the date format from the row above so we reuse the existing style
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A8").Value = 44582
$ws.Range("B8").Value = "Random Forest"
$ws.Range("C8").Value = "Tuned"
$ws.Range("D8").Value = 9.02

# Update selection to match the diff (A9 selected)
$ws.Range("A9").Select()
